# Adds 7 new Delta Smelt catch records (rows 256-262) to the
# "Delta Smelt Catch Data" worksheet, matching the diet-cages /
# recapture data collected for 24-33 / 24-35 releases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Delta Smelt Catch Data")

# xlCenter / xlLeft alignment constants used below
$xlCenter = -4108
$xlLeft = -4131

$rows = @(
    @{ r = 256; A = 45362; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-33-SM02";  E = "Suisun Marsh";                  F = 1113;  G = 13.5; H = 70.2; I = 47.32; J = 72; K = "Adult"; L = 38.146929999999998; M = -122.05933;              N = "Directed Outflow Project"; O = "S261"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Left/Orange/Posterior Dorsal";  S = "BY2023 3a"; T = 45315; U = "Rio Vista"; V = "Hard (carboy)" },
    @{ r = 257; A = 45362; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-33-SSC03"; E = "Lower Sac River Ship Channel";  F = 334.2; G = 12.6; H = 55.9; I = 34.11; J = 66; K = "Adult"; L = 38.289020000000001; M = -121.65742;              N = "Directed Outflow Project"; O = "S332"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Right/Blue/Anterior Dorsal";   S = "BY2023 1b"; T = 45274; U = "Rio Vista"; V = "Soft (carboy)" },
    @{ r = 258; A = 45364; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-33-LSR02"; E = "Lower Sacramento River";        F = 195.3; G = 12.1; H = 43.5; I = 22.53; J = 79; K = "Adult"; L = 38.064660000000003; M = -121.79474999999999;     N = "Directed Outflow Project"; O = "S333"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Right/Orange/Anterior Dorsal"; S = "BY2023 3b"; T = 45316; U = "Rio Vista"; V = "Soft (carboy)" },
    @{ r = 259; A = 45376; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-35-LSR02"; E = "Lower Sacramento River";        F = 199.9; G = 14.2; H = 26.6; I = 15.37; J = 76; K = "Adult"; L = 38.06373;              M = -121.81019999999999;     N = "Directed Outflow Project"; O = "S334"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Left/Orange/Posterior Dorsal";  S = "BY2023 3a"; T = 45315; U = "Rio Vista"; V = "Hard (carboy)" },
    @{ r = 260; A = 45376; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-35-SM03";  E = "Suisun Marsh";                  F = 864;   G = 15.1; H = 45.3; I = 28.32; J = 85; K = "Adult"; L = 38.186950000000003; M = -121.97835000000001;     N = "Directed Outflow Project"; O = "S262"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Left/Orange/Posterior Dorsal";  S = "BY2023 3a"; T = 45315; U = "Rio Vista"; V = "Hard (carboy)" },
    @{ r = 261; A = 45376; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-35-SM03";  E = "Suisun Marsh";                  F = 880;   G = 15.1; H = 45.3; I = 27.63; J = 81; K = "Adult"; L = 38.186819999999997; M = -121.97931;              N = "Directed Outflow Project"; O = "S263"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Left/Orange/Posterior Dorsal";  S = "BY2023 3a"; T = 45315; U = "Rio Vista"; V = "Hard (carboy)" },
    @{ r = 262; A = 45377; B = "USFWS EDSM"; C = "Kodiak Trawl"; D = "24-35-SSC03"; E = "Lower Sac River Ship Channel";  F = 630;   G = 15;   H = 30.5; I = 18.87; J = 61; K = "Adult"; L = 38.314129999999999; M = -121.65241;              N = "Directed Outflow Project"; O = "S335"; P = "Liquid nitrogen"; Q = "UC Davis Aquatic Health Lab"; R = "VIE-Right/Blue/Anterior Dorsal";   S = "BY2023 1b"; T = 45274; U = "Rio Vista"; V = "Soft (carboy)" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    # column W (23) intentionally left blank for these rows

    # --- formatting, matching the look of the pasted-in source rows ---
    # Whole row: black font (pasted data keeps explicit black color rather
    # than the theme color used by the rest of the sheet)
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 31))
    $rowRange.Font.Color = 0

    # Most columns are centered
    $centerRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 23))
    $centerRange.HorizontalAlignment = $xlCenter

    # Date columns (SampleDate / ReleaseDate) use the short-date format
    $ws.Cells.Item($r, 1).NumberFormat = "mm-dd-yy"
    $ws.Cells.Item($r, 20).NumberFormat = "mm-dd-yy"

    # SpecialStudy column wraps text
    $ws.Cells.Item($r, 14).WrapText = $true

    # Trailing blank helper columns (X:AE) - X left aligned, Y:AE unformatted
    $ws.Cells.Item($r, 24).HorizontalAlignment = $xlLeft
    $tailRange = $ws.Range($ws.Cells.Item($r, 25), $ws.Cells.Item($r, 31))
    $tailRange.Font.Color = 0
}

# Update the active selection to reflect where the user ended up after
# entering the new data.
$ws.Activate()
$ws.Range("H270").Select()
